# Update the "consensus_estimates" sheet with refreshed quarterly/FY data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("consensus_estimates")

# Data rows (period, year, quarter, revenue, ebitda, eps, analysts_n) for rows 2..13.
# Empty numeric slots are represented with $null so the cell stays blank (but numeric-typed).
$data = @(
    @("Q3-2025", 2025, 3,    1300.2,  16.3,              -63.1,              -298.2),
    @("Q2-2025", 2025, 2,    1397.9,  180.8,              103,                -4.7),
    @("Q1-2025", 2025, 1,    862.5,   -66.59999999999999, -142.7,             -202.1),
    @("Q4-2024", 2024, 4,    1306.4,  155.4,              77.09999999999999, -135.6),
    @("Q3-2024", 2024, 3,    1348.8,  152.7,              71.90000000000001, -20.7),
    @("Q2-2024", 2024, 2,    $null,   $null,              $null,              $null),
    @("Q1-2024", 2024, 1,    $null,   $null,              $null,              $null),
    @("FY-2024",  2024, $null, 4637.2, 312.7,             -6.8,               -352.6),
    @("FY-2023",  2023, $null, 4812.6, 404.1,             39.1,               -396.6),
    @("FY-2022",  2022, $null, 3911.4, 34.9,              -361.1,             -973.6),
    @("FY-2021",  2021, $null, 2527.9, -321.3,            -746.3,             -1269.1),
    @("FY-2020",  2020, $null, $null,  $null,             $null,              $null)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $ws.Cells.Item($rowIndex, 7).Value = $row[6]
    $rowIndex++
}
